# daily auto push: 2026-02-12 03:20 UTC
# Insert a new data row for 2026/02/12 07:00 (rank 201) right before the
# existing "2026/12/29" block, shifting all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 809 (and everything below it) down by one row.
$ws.Rows.Item(809).Insert()

# Column A holds plain text dates (e.g. "2026/02/12"), not real date
# serials, in this workbook. Force text formatting before assigning the
# value so Excel doesn't auto-convert the string into a date number, then
# copy the neighboring cell's (unformatted) style back so no stray
# number-format style sticks around on the new cell.
$ws.Range("A809").NumberFormat = "@"
$ws.Range("A809").Value = "2026/02/12"
$ws.Range("A809").Style = $ws.Range("A808").Style

$ws.Range("B809").Value = "木"
$ws.Range("C809").Value = 7
$ws.Range("D809").Value = 201
